$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$updates = @{
    2  = -3
    3  = -5
    10 = 4
    15 = 5
    19 = -1
    21 = -2
    39 = -1
    42 = 5
    43 = 2
    49 = 0
    50 = 0
    56 = 1
    59 = 2
    68 = 5
    69 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
